# Add a third column ("test") of data used for multidimensional linear
# regression, alongside the existing "temp" / "wind_speed" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("C1").Value = "test"

# Data rows 2-31
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 17
$ws.Range("C5").Value = 12
$ws.Range("C6").Value = 15
$ws.Range("C7").Value = 7
$ws.Range("C8").Value = 9
$ws.Range("C9").Value = 12
$ws.Range("C10").Value = 15
$ws.Range("C11").Value = 17
$ws.Range("C12").Value = 19.6666666666667
$ws.Range("C13").Value = 22.1666666666667
$ws.Range("C14").Value = 24.6666666666667
$ws.Range("C15").Value = 27.1666666666667
$ws.Range("C16").Value = 29.6666666666667
$ws.Range("C17").Value = 22
$ws.Range("C18").Value = 23
$ws.Range("C19").Value = 24
$ws.Range("C20").Value = 25
$ws.Range("C21").Value = 26
$ws.Range("C22").Value = 27
$ws.Range("C23").Value = 19
$ws.Range("C24").Value = 12
$ws.Range("C25").Value = 5
$ws.Range("C26").Value = 8
$ws.Range("C27").Value = 7
$ws.Range("C28").Value = 10
$ws.Range("C29").Value = 15
$ws.Range("C30").Value = "21``"
$ws.Range("C31").Value = 12

# Update selection to match the final cursor position recorded in the file
$ws.Range("C31").Select()
